$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.996.38'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '3.522.08'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = "'592.43"
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').Value = "'133.81"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.17%  '

$ws.Range('D7').Value = '3.521.02'
$ws.Range('E7').Value = '  -0.82%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('E9').Value = '  -0.57%  '

$ws.Range('E10').Value = '  +1.97%  '

$ws.Range('D11').Value = "'7.13"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.25%  '

$ws.Range('E12').Value = '  +0.01%  '

$ws.Range('D13').Value = '4.123.33'
$ws.Range('E13').Value = '  -0.79%  '

$ws.Range('D14').Value = "'27.67"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.71%  '

$ws.Range('D15').Value = "'0.0000181"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.05%  '

$ws.Range('E16').Value = '  +0.51%  '

$ws.Range('D17').Value = '3.526.41'
$ws.Range('E17').Value = '  -0.90%  '

$ws.Range('D18').Value = '64.980.80'
$ws.Range('E18').Value = '  +0.41%  '

$ws.Range('D19').Value = "'10.13"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.14%  '

$ws.Range('D20').Value = "'14.40"
$ws.Range('D20').Style = 'Normal'

$ws.Range('E21').Value = '  -2.27%  '

$ws.Range('D22').Value = "'391.80"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.33%  '

$ws.Range('D23').Value = "'0.580"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.89%  '

$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.667.47'
$ws.Range('E24').Value = '  -0.73%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'74.79"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.60%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('E27').Value = '  -3.34%  '

$ws.Range('E28').Value = '  +1.73%  '

$ws.Range('D29').Value = "'1.60"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.09%  '

$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('E31').Value = '  -0.74%  '

$ws.Range('E32').Value = '  -0.17%  '

$ws.Range('D33').Value = '3.531.26'
$ws.Range('E33').Value = '  -0.81%  '

$ws.Range('D34').Value = "'24.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.73%  '

$ws.Range('E36').Value = '  +0.04%  '

$ws.Range('E37').Value = '  +6.65%  '

$ws.Range('E38').Value = '  +2.86%  '

$ws.Range('D39').Value = "'6.97"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.80%  '

$ws.Range('D40').Value = "'168.37"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.81%  '

$ws.Range('D41').Value = "'0.0810"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.87%  '

$ws.Range('D42').Value = "'0.823"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.41%  '

$ws.Range('E43').Value = '  +5.89%  '

$ws.Range('D44').Value = "'25.85"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.01%  '

$ws.Range('D45').Value = "'42.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.50%  '

$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').Value = "'4.43"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.34%  '

$ws.Range('E48').Value = '  +1.11%  '

$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('D50').Value = '2.426.00'
$ws.Range('E50').Value = '  -1.00%  '

$ws.Range('D51').Value = "'0.910"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.58%  '
